# Update the "Naplánované hodiny" (planned hours) formulas for each member
# with the new contributing values, matching the newly recalculated totals.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

$ws.Range("B5").Formula  = "=5+12+14+15+10+19.5+3+2"
$ws.Range("B6").Formula  = "=15.5+1+11+5+6+10+15+2+3+2+2"
$ws.Range("B7").Formula  = "=18.5+13+16+14+5+3+2"
$ws.Range("B8").Formula  = "=14.5+13+6+15+5+8+3+2+11"
$ws.Range("B9").Formula  = "=14.75+5+16+6+15+2+3+2+13+4"
$ws.Range("B10").Formula = "=14.5+7+13+9+10+5+3+2+13"

# Move/update the active selection to reflect where the user last clicked.
[void]$ws.Range("C12").Select()
